$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.276052666666667
$ws.Range("H2").Value = 6.828158
$ws.Range("I2").Value = 0.005247614157263819
$ws.Range("J2").Value = 0.005247614157263819
$ws.Range("M2").Value = 16.27546433333333
$ws.Range("N2").Value = 48.826393
$ws.Range("O2").Value = 0.06628560529319844
$ws.Range("P2").Value = 0.06628560529319844
$ws.Range("Q2").Value = 37.04381399712155
$ws.Range("R2").Value = 333.394325974094
$ws.Range("S2").Value = 0.0003478412807593896
$ws.Range("T2").Value = 0.0003478412807593897

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.276052666666667
$ws.Range("H3").Value = 6.828158
$ws.Range("I3").Value = 0.005247614157263819
$ws.Range("J3").Value = 0.005247614157263819
$ws.Range("O3").Value = 0.3480686258826592
$ws.Range("P3").Value = 0.3480686258826592
$ws.Range("Q3").Value = 194.5186949473922
$ws.Range("R3").Value = 1750.66825452653
$ws.Range("S3").Value = 0.001826529848881206
$ws.Range("T3").Value = 0.001826529848881206

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.276052666666667
$ws.Range("H4").Value = 6.828158
$ws.Range("I4").Value = 0.005247614157263819
$ws.Range("J4").Value = 0.005247614157263819
$ws.Range("M4").Value = 42.61351133333333
$ws.Range("N4").Value = 127.840534
$ws.Range("O4").Value = 0.17355341356458
$ws.Range("P4").Value = 0.17355341356458
$ws.Range("Q4").Value = 96.99059610626355
$ws.Range("R4").Value = 872.9153649563719
$ws.Range("S4").Value = 0.0009107413500629523
$ws.Range("T4").Value = 0.0009107413500629527

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.276052666666667
$ws.Range("H5").Value = 6.828158
$ws.Range("I5").Value = 0.005247614157263819
$ws.Range("J5").Value = 0.005247614157263819
$ws.Range("M5").Value = 101.183272
$ws.Range("N5").Value = 303.549816
$ws.Range("O5").Value = 0.4120923552595624
$ws.Range("P5").Value = 0.4120923552595624
$ws.Range("Q5").Value = 230.2984560576587
$ws.Range("R5").Value = 2072.686104518928
$ws.Range("S5").Value = 0.002162501677560271
$ws.Range("T5").Value = 0.002162501677560271

$ws.Range("I6").Value = 0.1062533062835484
$ws.Range("J6").Value = 0.1062533062835484
$ws.Range("M6").Value = 16.27546433333333
$ws.Range("N6").Value = 48.826393
$ws.Range("O6").Value = 0.06628560529319844
$ws.Range("P6").Value = 0.06628560529319844
$ws.Range("Q6").Value = 750.0604268129451
$ws.Range("R6").Value = 6750.543841316506
$ws.Range("S6").Value = 0.007043064721408609
$ws.Range("T6").Value = 0.007043064721408611

$ws.Range("I7").Value = 0.1062533062835484
$ws.Range("J7").Value = 0.1062533062835484
$ws.Range("O7").Value = 0.3480686258826592
$ws.Range("P7").Value = 0.3480686258826592
$ws.Range("S7").Value = 0.036983442313604
$ws.Range("T7").Value = 0.03698344231360401

$ws.Range("I8").Value = 0.1062533062835484
$ws.Range("J8").Value = 0.1062533062835484
$ws.Range("M8").Value = 42.61351133333333
$ws.Range("N8").Value = 127.840534
$ws.Range("O8").Value = 0.17355341356458
$ws.Range("P8").Value = 0.17355341356458
$ws.Range("Q8").Value = 1963.858470889603
$ws.Range("R8").Value = 17674.72623800643
$ws.Range("S8").Value = 0.01844062400803266
$ws.Range("T8").Value = 0.01844062400803266

$ws.Range("I9").Value = 0.1062533062835484
$ws.Range("J9").Value = 0.1062533062835484
$ws.Range("M9").Value = 101.183272
$ws.Range("N9").Value = 303.549816
$ws.Range("O9").Value = 0.4120923552595624
$ws.Range("P9").Value = 0.4120923552595624
$ws.Range("Q9").Value = 4663.066234443141
$ws.Range("R9").Value = 41967.59610998827
$ws.Range("S9").Value = 0.04378617524050311
$ws.Range("T9").Value = 0.04378617524050312

$ws.Range("G10").Value = 41.187613
$ws.Range("H10").Value = 123.562839
$ws.Range("I10").Value = 0.09496120377532416
$ws.Range("J10").Value = 0.09496120377532417
$ws.Range("M10").Value = 16.27546433333333
$ws.Range("N10").Value = 48.826393
$ws.Range("O10").Value = 0.06628560529319844
$ws.Range("P10").Value = 0.06628560529319844
$ws.Range("Q10").Value = 670.3475263566363
$ws.Range("R10").Value = 6033.127737209727
$ws.Range("S10").Value = 0.006294560871618123
$ws.Range("T10").Value = 0.006294560871618123

$ws.Range("G11").Value = 41.187613
$ws.Range("H11").Value = 123.562839
$ws.Range("I11").Value = 0.09496120377532416
$ws.Range("J11").Value = 0.09496120377532417
$ws.Range("O11").Value = 0.3480686258826592
$ws.Range("P11").Value = 0.3480686258826592
$ws.Range("Q11").Value = 3520.024314943319
$ws.Range("R11").Value = 31680.21883448987
$ws.Range("S11").Value = 0.03305301571024027
$ws.Range("T11").Value = 0.03305301571024027

$ws.Range("G12").Value = 41.187613
$ws.Range("H12").Value = 123.562839
$ws.Range("I12").Value = 0.09496120377532416
$ws.Range("J12").Value = 0.09496120377532417
$ws.Range("M12").Value = 42.61351133333333
$ws.Range("N12").Value = 127.840534
$ws.Range("O12").Value = 0.17355341356458
$ws.Range("P12").Value = 0.17355341356458
$ws.Range("Q12").Value = 1755.148813368447
$ws.Range("R12").Value = 15796.33932031602
$ws.Range("S12").Value = 0.01648084107140919
$ws.Range("T12").Value = 0.01648084107140919

$ws.Range("G13").Value = 41.187613
$ws.Range("H13").Value = 123.562839
$ws.Range("I13").Value = 0.09496120377532416
$ws.Range("J13").Value = 0.09496120377532417
$ws.Range("M13").Value = 101.183272
$ws.Range("N13").Value = 303.549816
$ws.Range("O13").Value = 0.4120923552595624
$ws.Range("P13").Value = 0.4120923552595624
$ws.Range("Q13").Value = 4167.497449209735
$ws.Range("R13").Value = 37507.47704288762
$ws.Range("S13").Value = 0.03913278612205658
$ws.Range("T13").Value = 0.03913278612205659

$ws.Range("G14").Value = 344.1819356666667
$ws.Range("H14").Value = 1032.545807
$ws.Range("I14").Value = 0.7935378757838636
$ws.Range("J14").Value = 0.7935378757838637
$ws.Range("M14").Value = 16.27546433333333
$ws.Range("N14").Value = 48.826393
$ws.Range("O14").Value = 0.06628560529319844
$ws.Range("P14").Value = 0.06628560529319844
$ws.Range("Q14").Value = 5601.720818120461
$ws.Range("R14").Value = 50415.48736308415
$ws.Range("S14").Value = 0.05260013841941231
$ws.Range("T14").Value = 0.05260013841941232

$ws.Range("G15").Value = 344.1819356666667
$ws.Range("H15").Value = 1032.545807
$ws.Range("I15").Value = 0.7935378757838636
$ws.Range("J15").Value = 0.7935378757838637
$ws.Range("O15").Value = 0.3480686258826592
$ws.Range("P15").Value = 0.3480686258826592
$ws.Range("Q15").Value = 29414.8821469922
$ws.Range("R15").Value = 264733.9393229298
$ws.Range("S15").Value = 0.2762056380099337
$ws.Range("T15").Value = 0.2762056380099338

$ws.Range("G16").Value = 344.1819356666667
$ws.Range("H16").Value = 1032.545807
$ws.Range("I16").Value = 0.7935378757838636
$ws.Range("J16").Value = 0.7935378757838637
$ws.Range("M16").Value = 42.61351133333333
$ws.Range("N16").Value = 127.840534
$ws.Range("O16").Value = 0.17355341356458
$ws.Range("P16").Value = 0.17355341356458
$ws.Range("Q16").Value = 14666.8008162601
$ws.Range("R16").Value = 132001.2073463409
$ws.Range("S16").Value = 0.1377212071350752
$ws.Range("T16").Value = 0.1377212071350752

$ws.Range("G17").Value = 344.1819356666667
$ws.Range("H17").Value = 1032.545807
$ws.Range("I17").Value = 0.7935378757838636
$ws.Range("J17").Value = 0.7935378757838637
$ws.Range("M17").Value = 101.183272
$ws.Range("N17").Value = 303.549816
$ws.Range("O17").Value = 0.4120923552595624
$ws.Range("P17").Value = 0.4120923552595624
$ws.Range("Q17").Value = 34825.45441404683
$ws.Range("R17").Value = 313429.0897264215
$ws.Range("S17").Value = 0.3270108922194424
$ws.Range("T17").Value = 0.3270108922194425
